# Insert a new data row at row 215 (this shifts existing rows 215..324 down to 216..325)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("215:215").Insert()

# Populate the newly inserted row 215 with its values
$ws.Range("A215").Value = 8
$ws.Range("B215").Value = "Terminal La Palmera de La Serena"
$ws.Range("C215").Value = "Coquimbo"
$ws.Range("D215").Value = 44846
$ws.Range("E215").Value = 4
$ws.Range("F215").Value = 100112012
$ws.Range("G215").Value = "Espinaca"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 1200
$ws.Range("K215").Value = 450
$ws.Range("L215").Value = 500
$ws.Range("M215").Value = 475
$ws.Range("N215").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O215").Value = "Provincia del Elquí"
$ws.Range("P215").Value = 950
$ws.Range("Q215").Value = 0.5
$ws.Range("R215").Value = "Hortaliza"
